{"js": "const replacements = [\n  [\"455\u00f77=\", \"693\u00f73=\"],\n  [\"380\u00f79=\", \"964\u00f77=\"],\n  [\"634\u00f78=\", \"118\u00f79=\"],\n  [\"246\u00f79=\", \"581\u00f75=\"],\n  [\"290\u00f79=\", \"912\u00f72=\"],\n  [\"264\u00f74=\", \"258\u00f74=\"],\n  [\"643\u00f77=\", \"682\u00f75=\"],\n  [\"840\u00f74=\", \"382\u00f73=\"],\n  [\"195\u00f73=\", \"850\u00f77=\"],\n  [\"313\u00f73=\", \"808\u00f76=\"],\n  [\"601\u00f75=\", \"394\u00f77=\"],\n  [\"370\u00f76=\", \"807\u00f73=\"],\n  [\"858\u00f75=\", \"671\u00f73=\"],\n  [\"407\u00f78=\", \"720\u00f72=\"],\n  [\"316\u00f79=\", \"487\u00f79=\"],\n  [\"473\u00f76=\", \"853\u00f76=\"],\n  [\"544\u00f72=\", \"545\u00f77=\"],\n  [\"568\u00f76=\", \"933\u00f72=\"],\n  [\"930\u00f72=\", \"337\u00f72=\"],\n  [\"467\u00f76=\", \"991\u00f75=\"],\n  [\"440\u00f75=\", \"408\u00f75=\"],\n  [\"607\u00f72=\", \"647\u00f78=\"],\n  [\"636\u00f74=\", \"116\u00f72=\"],\n  [\"417\u00f78=\", \"925\u00f73=\"],\n  [\"631\u00f77=\", \"674\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('455\u00f77=', '693\u00f73='),\n    @('380\u00f79=', '964\u00f77='),\n    @('634\u00f78=', '118\u00f79='),\n    @('246\u00f79=', '581\u00f75='),\n    @('290\u00f79=', '912\u00f72='),\n    @('264\u00f74=', '258\u00f74='),\n    @('643\u00f77=', '682\u00f75='),\n    @('840\u00f74=', '382\u00f73='),\n    @('195\u00f73=', '850\u00f77='),\n    @('313\u00f73=', '808\u00f76='),\n    @('601\u00f75=', '394\u00f77='),\n    @('370\u00f76=', '807\u00f73='),\n    @('858\u00f75=', '671\u00f73='),\n    @('407\u00f78=', '720\u00f72='),\n    @('316\u00f79=', '487\u00f79='),\n    @('473\u00f76=', '853\u00f76='),\n    @('544\u00f72=', '545\u00f77='),\n    @('568\u00f76=', '933\u00f72='),\n    @('930\u00f72=', '337\u00f72='),\n    @('467\u00f76=', '991\u00f75='),\n    @('440\u00f75=', '408\u00f75='),\n    @('607\u00f72=', '647\u00f78='),\n    @('636\u00f74=', '116\u00f72='),\n    @('417\u00f78=', '925\u00f73='),\n    @('631\u00f77=', '674\u00f79='),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
